$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update aggregate metrics after trade #26 closed
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.71   # Current Capital
$summary.Range("B4").Value = 0.7       # Total P&L $
$summary.Range("B5").Value = 0.54      # Total P&L %
$summary.Range("B6").Value = 26        # Total Trades
$summary.Range("B7").Value = 11        # Winning Trades
$summary.Range("B9").Value = 42.31     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.71     # Capital
$status.Range("D4").Value = 26         # Trades
$status.Range("E4").Value = 0.7        # P&L $
$status.Range("F4").Value = 0.71       # P&L %
$status.Range("G4").Value = 42.31      # Win Rate %

# ---------------------------------------------------------------------------
# Append the newly closed trade (#26) to both "All Trades" and
# "MarketMaking" sheets, which keep an identical trade log.
# ---------------------------------------------------------------------------
$newRow = @{
    A = 26
    B = "2026-02-17"
    C = "12:37:37"
    D = "MarketMaking"
    E = "DOWN"
    F = 0.33
    G = 0.341461
    H = "CLOSED"
    I = 3.4732
    J = 0.01
    K = 100.71
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.13
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A27").Value = $newRow.A

    # "2026-02-17" looks like a date, and a plain .Value assignment would
    # get auto-converted by Excel into a date serial number. Force it to be
    # stored as text (matching the original sheet, where Date is a plain
    # string), then reset the cell style back to Normal so no stray
    # number-format style lingers on the cell.
    $ws.Range("B27").NumberFormat = "@"
    $ws.Range("B27").Value = $newRow.B
    $ws.Range("B27").Style = "Normal"

    $ws.Range("C27").Value = $newRow.C
    $ws.Range("D27").Value = $newRow.D
    $ws.Range("E27").Value = $newRow.E
    $ws.Range("F27").Value = $newRow.F
    $ws.Range("G27").Value = $newRow.G
    $ws.Range("H27").Value = $newRow.H
    $ws.Range("I27").Value = $newRow.I
    $ws.Range("J27").Value = $newRow.J
    $ws.Range("K27").Value = $newRow.K
    $ws.Range("L27").Value = $newRow.L
    $ws.Range("M27").Value = $newRow.M
    $ws.Range("N27").Value = $newRow.N
    $ws.Range("O27").Value = $newRow.O
    $ws.Range("P27").Value = $newRow.P
    $ws.Range("Q27").Value = $newRow.Q
}
